# Update odds values on row 3 of the active worksheet to reflect
# refreshed FlashScore data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "H3"  = 4.33
    "I3"  = 5
    "J3"  = 2.1
    "K3"  = 2.5
    "S3"  = 1.29
    "T3"  = 3.5
    "Z3"  = 12
    "AC3" = 17
    "AD3" = 8.5
    "AI3" = 29
    "AT3" = 3.5
    "AU3" = 7.5
    "AV3" = 41
    "AX3" = 26
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
